$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Bottom-up row insertions so earlier row numbers stay stable while we work ---

# 1) After old row 41 (Sreg_wrapper.c / double B / double B=0) insert two new rows
#    for "double c;"/"double c=0;" and "double d;"/"double d=0;"
$ws.Rows.Item(42).Insert() | Out-Null
$ws.Rows.Item(42).Insert() | Out-Null

$ws.Cells.Item(42,1).Value = "Sreg_wrapper.c"
$ws.Cells.Item(42,2).Value = "double c;"
$ws.Cells.Item(42,3).Value = "double c=0;"

$ws.Cells.Item(43,1).Value = "Sreg_wrapper.c"
$ws.Cells.Item(43,2).Value = "double d;"
$ws.Cells.Item(43,3).Value = "double d=0;"

# 2) After old row 35 (Sreg.c / double A / double A=0) insert one new row
#    for "double B;"/"double B=0;"
$ws.Rows.Item(36).Insert() | Out-Null

$ws.Cells.Item(36,1).Value = "Sreg.c"
$ws.Cells.Item(36,2).Value = "double B;"
$ws.Cells.Item(36,3).Value = "double B=0;"

# 3) After old row 30 (MMregcore.c / double A2 / double A2=0) insert two new rows
#    for "double B2;"/"double B2=0;" and "double ctun;"/"double ctun=0;"
$ws.Rows.Item(31).Insert() | Out-Null
$ws.Rows.Item(31).Insert() | Out-Null

$ws.Cells.Item(31,1).Value = "MMregcore.c"
$ws.Cells.Item(31,2).Value = "double B2;"
$ws.Cells.Item(31,3).Value = "double B2=0;"

$ws.Cells.Item(32,1).Value = "MMregcore.c"
$ws.Cells.Item(32,2).Value = "double ctun;"
$ws.Cells.Item(32,3).Value = "double ctun=0;"

# --- In-place updates (row numbers below are final, post-insertion, numbers) ---

# Row 29: double a2 -> double a2; / double a2=0 -> double a2=0;
$ws.Cells.Item(29,2).Value = "double a2;"
$ws.Cells.Item(29,3).Value = "double a2=0;"

# Row 30: double A2 -> double A2; / double A2=0 -> double A2=0;
$ws.Cells.Item(30,2).Value = "double A2;"
$ws.Cells.Item(30,3).Value = "double A2=0;"

# Row 37 (was 35 before insertions): double A -> double A; / double A=0 -> double A=0;
$ws.Cells.Item(37,2).Value = "double A;"
$ws.Cells.Item(37,3).Value = "double A=0;"

# Row 40 (was 37 before insertions): double d -> double d; / double d=0 -> double d=0;
$ws.Cells.Item(40,2).Value = "double d;"
$ws.Cells.Item(40,3).Value = "double d=0;"

# Row 43 (was 40 before insertions): double A -> double A; / double A=0 -> double A=0;
$ws.Cells.Item(43,2).Value = "double A;"
$ws.Cells.Item(43,3).Value = "double A=0;"

# Row 44 (was 41 before insertions): double B -> double B; / double B=0 -> double B=0;
$ws.Cells.Item(44,2).Value = "double B;"
$ws.Cells.Item(44,3).Value = "double B=0;"

# --- Sheet view bookkeeping to match the saved workbook state ---
$ws.Range("B45").Select() | Out-Null
[void]($excel.ActiveWindow.ScrollRow = 28)
